$d = $word.ActiveDocument

# Locate the "Edison Achalma" paragraph that carries the "Author" style
# (the byline right under the "Editar: Editar" title heading — not the
# later mentions of the same name in the body text further down).
$range = $d.Content
$found = $range.Find.Execute("Edison Achalma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Move to the end of that match (still inside the "Author" paragraph,
# right before its paragraph mark) and split off a brand-new paragraph
# right after it by inserting a paragraph mark.
$range.Collapse(0)
$range.InsertAfter("`r")
$range.Collapse(0)

# $range now sits at the very start of the freshly inserted (empty)
# paragraph, which inherited the "Author" style from its predecessor.
# Re-seat on a fresh Range at that exact position and fill it in with
# the affiliation line via InsertXML so the run is emitted exactly like
# the rest of the document (a single run with xml:space="preserve").
$insertAt = $d.Range($range.Start, $range.Start + 1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga</w:t></w:r></w:p>'
$insertAt.InsertXML($xml)
